$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "35.418.21"
$ws.Range("E2").Value = "  +1.23%  "
$ws.Range("D3").Value = "1.872.31"
$ws.Range("E3").Value = "  +1.62%  "
$ws.Range("E4").Value = "  +0.55%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "241.11"
$ws.Range("E5").Value = "  +4.13%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.631"
$ws.Range("E6").Value = "  +1.95%  "
$ws.Range("E7").Value = "  +0.40%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "42.50"
$ws.Range("E8").Value = "  +7.17%  "
$ws.Range("E9").Value = "  +0.67%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0704"
$ws.Range("E10").Value = "  +2.11%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0989"
$ws.Range("E11").Value = "  +0.39%  "
$ws.Range("D12").Value = "2.142.03"
$ws.Range("E12").Value = "  +1.55%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.897.39"
$ws.Range("E13").Value = "  +2.93%  "
$ws.Range("B14").Value = "Chainlink"
$ws.Range("C14").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "11.62"
$ws.Range("E14").Value = "  +1.14%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.684"
$ws.Range("E15").Value = "  +1.72%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.76"
$ws.Range("E16").Value = "  +2.47%  "
$ws.Range("D17").Value = "35.412.25"
$ws.Range("E17").Value = "  +1.21%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "70.99"
$ws.Range("E18").Value = "  +1.63%  "
$ws.Range("D19").Value = "0.0₃0804"
$ws.Range("E19").Value = "  +1.86%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "242.83"
$ws.Range("E20").Value = "  +1.29%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.37"
$ws.Range("E21").Value = "  +1.41%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.78"
$ws.Range("E22").Value = "  +2.11%  "
$ws.Range("E23").Value = "  +0.39%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.28"
$ws.Range("E24").Value = "  -0.42%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "170.38"
$ws.Range("E25").Value = "  -1.12%  "
$ws.Range("B26").Value = "PancakeSwap"
$ws.Range("C26").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.95"
$ws.Range("E26").Value = "  +27.17%  "
$ws.Range("B27").Value = "Cosmos"
$ws.Range("C27").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.22"
$ws.Range("E27").Value = "  +5.07%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.77"
$ws.Range("E28").Value = "  +1.69%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.126"
$ws.Range("E29").Value = "  +1.51%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0565"
$ws.Range("E30").Value = "  +2.47%  "
$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.06"
$ws.Range("E31").Value = "  +2.85%  "
$ws.Range("B32").Value = "BinanceUSD"
$ws.Range("C32").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.01"
$ws.Range("E32").Value = "  +0.59%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.08"
$ws.Range("E33").Value = "  +2.74%  "
$ws.Range("B34").Value = "ImmutableX"
$ws.Range("C34").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.875"
$ws.Range("E34").Value = "  +22.22%  "
$ws.Range("B35").Value = "WEMIXToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.70"
$ws.Range("E35").Value = "  +5.70%  "
$ws.Range("E36").Value = "  +5.48%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.33"
$ws.Range("E37").Value = "  +9.06%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.10"
$ws.Range("E38").Value = "  +3.64%  "
$ws.Range("E39").Value = "  +4.72%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "91.04"
$ws.Range("E40").Value = "  +1.25%  "
$ws.Range("D41").Value = "1.359.54"
$ws.Range("E41").Value = "  +0.49%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "15.38"
$ws.Range("E42").Value = "  +4.55%  "
$ws.Range("B43").Value = "MultiversX"
$ws.Range("C43").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "49.73"
$ws.Range("E43").Value = "  +48.25%  "
$ws.Range("B44").Value = "Kaspa"
$ws.Range("C44").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0590"
$ws.Range("E44").Value = "  +11.31%  "
$ws.Range("E45").Value = "  +3.28%  "
$ws.Range("B46").Value = "Gas"
$ws.Range("C46").Value = "https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.59"
$ws.Range("E46").Value = "  +46.79%  "
$ws.Range("B47").Value = "HuobiToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.42"
$ws.Range("E47").Value = "  +0.07%  "
$ws.Range("B48").Value = "FraxShare"
$ws.Range("C48").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "6.72"
$ws.Range("E48").Value = "  +7.99%  "
$ws.Range("B49").Value = "MXToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.72"
$ws.Range("E49").Value = "  -1.25%  "
$ws.Range("B50").Value = "RocketPoolETH"
$ws.Range("C50").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D50").Value = "2.059.43"
$ws.Range("E50").Value = "  +1.67%  "
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0689"
$ws.Range("E51").Value = "  +2.49%  "
